$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 19: fill in the spent-time value
$ws.Range("E19").Value = 1.3

# Row 30: "Jak starczy czasu" + chain-rule link (added first so the shared
# strings land in the same order as row 29's additions)
$ws.Range("F30").Value = "Jak starczy czasu"
$ws.Range("G30").Value = "https://pl.khanacademy.org/math/differential-calculus/dc-chain#dc-chain-rule"

# Row 29: local-extrema repetition note + link
$ws.Range("K29").Value = "https://pl.khanacademy.org/math/ap-calculus-ab/ab-diff-analytical-applications-new/ab-5-4/e/critical-numbers"
$ws.Range("J29").Value = "Powt z ekstremów lok"

# Move the active selection to H39
$ws.Range("H39").Select()
